$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 29 - this shifts the existing row 29..61 down to 30..62
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new weekly record
$ws.Cells.Item(29, 1).Value = 6
$ws.Cells.Item(29, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(29, 3).Value = "Metropolitana"
$ws.Cells.Item(29, 4).Value = 44494
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = "Fruta"
$ws.Cells.Item(29, 7).Value = 100108
$ws.Cells.Item(29, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(29, 9).Value = 100108007
$ws.Cells.Item(29, 10).Value = "Coco"
$ws.Cells.Item(29, 11).Value = "Sin especificar"
$ws.Cells.Item(29, 12).Value = "Primera"
$ws.Cells.Item(29, 13).Value = 150
$ws.Cells.Item(29, 14).Value = 24000
$ws.Cells.Item(29, 15).Value = 24000
$ws.Cells.Item(29, 16).Value = 24000
$ws.Cells.Item(29, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(29, 18).Value = "Perú"
$ws.Cells.Item(29, 19).Value = 1200
$ws.Cells.Item(29, 20).Value = 20
